# Updated cryptos list (price + 1h volume change) to match the latest
# scrape. Price cells that look like plain decimals (e.g. "0.517",
# "588.92") are written with a leading apostrophe so Excel stores them
# as literal text instead of coercing to a number, then the style is
# reset to "Normal" so no stray quote-prefix formatting is left on the
# cell (matches the original inline-string cells, which carry no `s`
# attribute).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.091.15"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").Value = "3.153.52"
$ws.Range("E3").Value = "  +1.40%  "

$ws.Range("D5").Value = "'588.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.48%  "

$ws.Range("D6").Value = "'138.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.90%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "3.151.10"
$ws.Range("E8").Value = "  +1.47%  "

$ws.Range("D9").Value = "'0.517"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.61%  "

$ws.Range("D10").Value = "'0.146"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.21%  "

$ws.Range("D11").Value = "'5.30"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.74%  "

$ws.Range("D12").Value = "'0.458"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.57%  "

$ws.Range("D13").Value = "'0.0000245"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.54%  "

$ws.Range("D14").Value = "'34.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.70%  "

$ws.Range("D15").Value = "3.666.56"
$ws.Range("E15").Value = "  +0.87%  "

$ws.Range("E16").Value = "  +0.67%  "

$ws.Range("D17").Value = "3.143.85"
$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("D18").Value = "63.019.19"
$ws.Range("E18").Value = "  -0.92%  "

$ws.Range("D19").Value = "'6.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.92%  "

$ws.Range("D20").Value = "'475.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.36%  "

$ws.Range("D21").Value = "'13.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.61%  "

$ws.Range("D22").Value = "'0.701"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.41%  "

$ws.Range("D23").Value = "'7.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.94%  "

$ws.Range("D24").Value = "'84.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.14%  "

$ws.Range("D25").Value = "'12.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.25%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "'2.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.13%  "

$ws.Range("D28").Value = "'7.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.86%  "

$ws.Range("D29").Value = "'7.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.16%  "

$ws.Range("E30").Value = "  +2.89%  "

$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").Value = "'26.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.00%  "

$ws.Range("E33").Value = "  -4.28%  "

$ws.Range("D34").Value = "'2.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.99%  "

$ws.Range("E35").Value = "  -2.85%  "

$ws.Range("D36").Value = "'52.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.10%  "

$ws.Range("D37").Value = "'5.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.04%  "

$ws.Range("D38").Value = "0.0₃0699"
$ws.Range("E38").Value = "  -5.16%  "

$ws.Range("D39").Value = "'0.0387"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.40%  "

$ws.Range("D40").Value = "'416.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.52%  "

$ws.Range("D41").Value = "'2.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.46%  "

$ws.Range("D42").Value = "'8.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.27%  "

$ws.Range("D43").Value = "2.925.12"
$ws.Range("E43").Value = "  +2.43%  "

$ws.Range("E44").Value = "  -6.58%  "

$ws.Range("D45").Value = "'0.259"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.21%  "

$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").Value = "'2.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.20%  "

$ws.Range("D48").Value = "'25.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.15%  "

$ws.Range("E49").Value = "  -0.28%  "

$ws.Range("D50").Value = "'2.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.50%  "

$ws.Range("D51").Value = "'120.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.10%  "
